# texts.xlsx - TouchGFX typography/translation sheet update
# "logica acc ok - logica autox da testare - schermata startup cambiata"
#
# Changes applied:
#  1. Typography sheet: row 10 ("Roboto-Black.ttf" / Bold) font Size 90 -> 100
#  2. Translation sheet: the "startup screen" text-id block (rows 256-265)
#     is re-shuffled / extended - a couple of new rows (284/285 "ready"
#     Bold entries) are added and existing rows shift so the Bold
#     ("<value>") pair now sits at rows 262 & 264, the "N"/"<value>" Extra
#     rows occupy 256-261, and row 263/265 carry the "READY" text.

$wb = $excel.ActiveWorkbook

# ---- 1. Typography sheet: D10 Size 90 -> 100 -------------------------
$typo = $wb.Worksheets.Item("Typography")
$typo.Range("D10").Value = 100

# ---- 2. Translation sheet: rows 256-265 ------------------------------
$tr = $wb.Worksheets.Item("Translation")

function Set-Row($r, $textId, $typography, $alignment, $value, $direction) {
    $tr.Range("B$r").Value = $textId
    $tr.Range("C$r").Value = $typography
    $tr.Range("D$r").Value = $alignment
    $tr.Range("E$r").Value = $value
    $tr.Range("F$r").Value = $direction
}

Set-Row 256 "SingleUseId277" "Extra" "Center" "<value>" "LTR"
Set-Row 257 "SingleUseId278" "Extra" "Left"   "N"       "LTR"
Set-Row 258 "SingleUseId279" "Extra" "Center" "<value>" "LTR"
Set-Row 259 "SingleUseId280" "Extra" "Left"   "N"       "LTR"
Set-Row 260 "SingleUseId281" "Extra" "Center" "<value>" "LTR"
Set-Row 261 "SingleUseId282" "Extra" "Left"   "N"       "LTR"
Set-Row 262 "SingleUseId276" "Bold"  "Center" "<value>" "LTR"
Set-Row 263 "SingleUseId283" "Bold"  "Left"   "READY"   "LTR"
Set-Row 264 "SingleUseId284" "Bold"  "Center" "<value>" "LTR"
Set-Row 265 "SingleUseId285" "Bold"  "Left"   "READY"   "LTR"
